# Updated data from Cannabis Commission
# - Append a new week's row of data to Sheet1
# - Add a new "underReview" worksheet with the latest (not-yet-approved) figures

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: append the new weekly row (row 7) ---------------------------
# Copy the formatting (date number format) from the cell above first, so the
# new date cell reuses the existing date style instead of creating a new one.
$ws1.Range("A6").Copy($ws1.Range("A7"))

$ws1.Cells.Item(7, 1).Value = 43256
$ws1.Cells.Item(7, 2).Value = 1002
$ws1.Cells.Item(7, 3).Value = 108
$ws1.Cells.Item(7, 4).Value = 61
$ws1.Cells.Item(7, 5).Value = 833
$ws1.Cells.Item(7, 6).Value = 1
$ws1.Cells.Item(7, 7).Value = 30
$ws1.Cells.Item(7, 8).Value = 22
$ws1.Cells.Item(7, 9).Value = 7
$ws1.Cells.Item(7, 10).Value = 14
$ws1.Cells.Item(7, 11).Value = 3
$ws1.Cells.Item(7, 12).Value = 30
$ws1.Cells.Item(7, 13).Value = 1

# Move the selection on Sheet1 to reflect where the author left off editing.
[void]$ws1.Range("F1:M1").Select()

# --- New worksheet: underReview ------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "underReview"

$ws2.Range("A1").Value = "date"
$ws2.Range("B1").Value = "appsReview"
$ws2.Range("C1").Value = "craftCooperative"
$ws2.Range("D1").Value = "cultivator"
$ws2.Range("E1").Value = "establishmentAgent"
$ws2.Range("F1").Value = "microbusiness"
$ws2.Range("G1").Value = "productManufacturer"
$ws2.Range("H1").Value = "researchFacility"
$ws2.Range("I1").Value = "retailer"
$ws2.Range("J1").Value = "transporter"

# Data row - reuse the date style from Sheet1 so it matches (m/d/yyyy).
$ws1.Range("A6").Copy($ws2.Range("A2"))
$ws2.Cells.Item(2, 1).Value = 43256
$ws2.Cells.Item(2, 2).Value = 51
$ws2.Cells.Item(2, 3).Value = 0
$ws2.Cells.Item(2, 4).Value = 18
$ws2.Cells.Item(2, 5).Value = "NA"
$ws2.Cells.Item(2, 6).Value = 2
$ws2.Cells.Item(2, 7).Value = 12
$ws2.Cells.Item(2, 8).Value = 3
$ws2.Cells.Item(2, 9).Value = 15
$ws2.Cells.Item(2, 10).Value = 1

# Column B is a little wider to fit the "appsReview" header.
$ws2.Columns.Item(2).ColumnWidth = 11.8

# Leave the cursor positioned under the header row, and make this the
# active/visible tab.
[void]$ws2.Range("A2").Select()
[void]$ws2.Activate()
